# Refresh market-price-derived columns (currentAveragePrice.., LevePrice.., LeveProfit..)
# for the Leve rows below, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row33 - "Glazed and Confused" (Clear Glass Lens)
$ws.Range("H33").Value = 331.05884
$ws.Range("I33").Value = 123.5
$ws.Range("K33").Value = 123.5
$ws.Range("M33").Value = 105.5

# ALC!row53 - "No Accounting for Waste" (Enchanted Electrum Ink)
$ws.Range("H53").Value = 280.1875
$ws.Range("I53").Value = 355.42856
$ws.Range("K53").Value = 355.42856
$ws.Range("M53").Value = 281.57144

# ALC!row76 - "Warding Off Temptation" (Enchanted Hardsilver Ink)
$ws.Range("H76").Value = 8000
$ws.Range("J76").Value = 8000
$ws.Range("L76").Value = 8000
$ws.Range("N76").Value = -8630

# ALC!row79 - "The Garden of Arcane Delights (L)" (Enchanted Hardsilver Ink)
$ws.Range("H79").Value = 8000
$ws.Range("J79").Value = 8000
$ws.Range("L79").Value = 8000
$ws.Range("N79").Value = -10184

# ALC!row116 - "Growing Up" (Growth Formula Kappa)
$ws.Range("H116").Value = 6393.1665
$ws.Range("J116").Value = 8996.333000000001
$ws.Range("L116").Value = 8996.333000000001
$ws.Range("N116").Value = -15880.333

# ALC!row121 - "Mindful Medicine" (Tincture of Mind)
$ws.Range("H121").Value = 548
$ws.Range("J121").Value = 548
$ws.Range("L121").Value = 1644
$ws.Range("N121").Value = -5138

# ALC!row132 - "Fast-forwarding Flora" (Growth Formula Lambda)
$ws.Range("H132").Value = 2355.3684
$ws.Range("I132").Value = 2322.0605
$ws.Range("J132").Value = 2575.2
$ws.Range("K132").Value = 6966.181500000001
$ws.Range("L132").Value = 7725.599999999999
$ws.Range("M132").Value = -4436.181500000001
$ws.Range("N132").Value = -12785.6

# ALC!row137 - "Cutting Edge of Culinary Quality" (Magnesia Whetstone)
$ws.Range("H137").Value = 12256.857
$ws.Range("I137").Value = 3949.75
$ws.Range("K137").Value = 11849.25
$ws.Range("M137").Value = -9299.25

# ALC!row138 - "All-night Crafting" (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 6874.7095
$ws.Range("J138").Value = 6670.324
$ws.Range("L138").Value = 20010.972
$ws.Range("N138").Value = -30290.972

$ws = $wb.Worksheets.Item("ARM")
# ARM!row2 - "Ain't Got No Ingots" (Bronze Ingot)
$ws.Range("H2").Value = 2122.7
$ws.Range("I2").Value = 2122.7
$ws.Range("K2").Value = 2122.7
$ws.Range("M2").Value = -2009.7

# ARM!row23 - "A Well-rounded Crew" (Iron Hoplon)
$ws.Range("H23").Value = 15333.333
$ws.Range("I23").Value = 15333.333
$ws.Range("K23").Value = 15333.333
$ws.Range("M23").Value = -15074.333

# ARM!row32 - "Ingot We Trust" (Steel Ingot)
$ws.Range("H32").Value = 24657.361
$ws.Range("I32").Value = 14883.588
$ws.Range("J32").Value = 30195.834
$ws.Range("K32").Value = 14883.588
$ws.Range("L32").Value = 30195.834
$ws.Range("M32").Value = -14596.588
$ws.Range("N32").Value = -30769.834

# ARM!row33 - "A Leg to Stand On" (Heavy Iron Flanchard)
$ws.Range("H33").Value = 1937.5
$ws.Range("I33").Value = 1937.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1937.5
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -1608.5

# ARM!row36 - "Hot for Teacher" (Heavy Iron Armor)
$ws.Range("H36").Value = 5508.5
$ws.Range("I36").Value = 5610.2
$ws.Range("K36").Value = 5610.2
$ws.Range("M36").Value = -5264.2

# ARM!row116 - "No Scope" (Titanbronze Ingot)
$ws.Range("H116").Value = 2122.7
$ws.Range("I116").Value = 2122.7
$ws.Range("K116").Value = 2122.7
$ws.Range("M116").Value = 171.3000000000002

# ARM!row132 - "Don't Bore Me, Ore Me" (Mountain Chromite Ingot)
$ws.Range("H132").Value = 3139
$ws.Range("I132").Value = 1666.4286
$ws.Range("K132").Value = 4999.2858
$ws.Range("M132").Value = -2469.2858

$ws = $wb.Worksheets.Item("BSM")
# BSM!row3 - "Hells Bells" (Bronze Ingot)
$ws.Range("H3").Value = 2122.7
$ws.Range("I3").Value = 2122.7
$ws.Range("K3").Value = 2122.7
$ws.Range("M3").Value = -2008.7

# BSM!row20 - "Smelt and Dealt" (Iron Ingot)
$ws.Range("H20").Value = 3863.5557
$ws.Range("I20").Value = 2388.5
$ws.Range("K20").Value = 2388.5
$ws.Range("M20").Value = -2141.5

# BSM!row64 - "With Bearings Straight" (Mythrite Nugget)
$ws.Range("H64").Value = 1609.1818
$ws.Range("I64").Value = 1361
$ws.Range("J64").Value = 1816
$ws.Range("K64").Value = 1361
$ws.Range("L64").Value = 1816
$ws.Range("M64").Value = -1136
$ws.Range("N64").Value = -2266

# BSM!row67 - "Bearing the Brunt (L)" (Mythrite Nugget)
$ws.Range("H67").Value = 1609.1818
$ws.Range("I67").Value = 1361
$ws.Range("J67").Value = 1816
$ws.Range("K67").Value = 1361
$ws.Range("L67").Value = 1816
$ws.Range("M67").Value = -581
$ws.Range("N67").Value = -3376

$ws = $wb.Worksheets.Item("CRP")
# CRP!row32 - "Daddy's Little Girl" (Viper-crested Round Shield)
$ws.Range("H32").Value = 2233.3333
$ws.Range("I32").Value = 2233.3333
$ws.Range("K32").Value = 2233.3333
$ws.Range("M32").Value = -1917.3333

# CRP!row60 - "Bowing to Greater Power" (Yew Longbow)
$ws.Range("H60").Value = 56999
$ws.Range("J60").Value = 69999
$ws.Range("L60").Value = 69999
$ws.Range("N60").Value = -71021

# CRP!row62 - "Splinter in the Sewers" (Cedar Lumber)
$ws.Range("H62").Value = 118285.57
$ws.Range("I62").Value = 8666.666999999999
$ws.Range("J62").Value = 200499.75
$ws.Range("K62").Value = 8666.666999999999
$ws.Range("L62").Value = 200499.75
$ws.Range("M62").Value = -8042.666999999999
$ws.Range("N62").Value = -201747.75

# CRP!row65 - "The Lumber of Their Discontent (L)" (Cedar Lumber)
$ws.Range("H65").Value = 118285.57
$ws.Range("I65").Value = 8666.666999999999
$ws.Range("J65").Value = 200499.75
$ws.Range("K65").Value = 43333.335
$ws.Range("L65").Value = 1002498.75
$ws.Range("M65").Value = -40213.335
$ws.Range("N65").Value = -1008738.75

# CRP!row94 - "Beech, Please" (Beech Lumber)
$ws.Range("H94").Value = 973
$ws.Range("I94").Value = 690
$ws.Range("J94").Value = 1086.2
$ws.Range("K94").Value = 690
$ws.Range("L94").Value = 1086.2
$ws.Range("M94").Value = -239
$ws.Range("N94").Value = -1988.2

# CRP!row105 - "Zelkova, My Love" (Zelkova Lumber)
$ws.Range("H105").Value = 4002.818
$ws.Range("I105").Value = 2068.3333
$ws.Range("K105").Value = 2068.3333
$ws.Range("M105").Value = -321.3332999999998

# CRP!row107 - "Built to Last" (White Oak Lumber)
$ws.Range("H107").Value = 798.86365
$ws.Range("I107").Value = 249.55556
$ws.Range("J107").Value = 1179.1538
$ws.Range("K107").Value = 249.55556
$ws.Range("L107").Value = 1179.1538
$ws.Range("M107").Value = 1670.44444
$ws.Range("N107").Value = -5019.1538

$ws = $wb.Worksheets.Item("GSM")
# GSM!row11 - "A Ringing Success" (Copper Ring)
$ws.Range("H11").Value = 8446224
$ws.Range("J11").Value = 8001804
$ws.Range("L11").Value = 8001804
$ws.Range("N11").Value = -8002082

# GSM!row21 - "Forever 21K" (Brass Ring)
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# GSM!row24 - "Bad Guys Eat Brass" (Brass Ring of Crafting)
$ws.Range("H24").Value = 32935.273
$ws.Range("I24").Value = 15000
$ws.Range("J24").Value = 34728.8
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 34728.8
$ws.Range("M24").Value = -14827
$ws.Range("N24").Value = -35074.8

# GSM!row28 - "You Burnt It, You Bought It" (Wind Brand)
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# GSM!row30 - "Dog Tags Are for Dogs" (Brass Ring)
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# GSM!row132 - "On Board for Lar" (Lar Ingot)
$ws.Range("H132").Value = 2936.5806
$ws.Range("I132").Value = 2035.44
$ws.Range("J132").Value = 6691.3335
$ws.Range("K132").Value = 6106.32
$ws.Range("L132").Value = 20074.0005
$ws.Range("M132").Value = -3576.32
$ws.Range("N132").Value = -25134.0005

$ws = $wb.Worksheets.Item("LTW")
# LTW!row100 - "Tiger in the Sack" (Tiger Leather)
$ws.Range("H100").Value = 2670.4
$ws.Range("J100").Value = 2000
$ws.Range("L100").Value = 2000
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("WVR")
# WVR!row100 - "Of Great Import" (Kudzu Thread)
$ws.Range("H100").Value = 2237.5
$ws.Range("I100").Value = 2640
$ws.Range("J100").Value = 1566.6666
$ws.Range("K100").Value = 5280
$ws.Range("L100").Value = 3133.3332
$ws.Range("M100").Value = -4739
$ws.Range("N100").Value = -4215.3332

# WVR!row126 - "A Polished Purchase" (Snow Linen)
$ws.Range("H126").Value = 118608.664
$ws.Range("I126").Value = 206496.6
$ws.Range("J126").Value = 8748.75
$ws.Range("K126").Value = 619489.8
$ws.Range("L126").Value = 26246.25
$ws.Range("M126").Value = -617019.8
$ws.Range("N126").Value = -31186.25
